$wb = $excel.ActiveWorkbook

# ================= Overview sheet =================
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Range("A5").Value = "d1e0eeed-5c24-49dd-bf01-14fb964d9451.md"
$wsOv.Range("C5").Value = ".md"
$wsOv.Range("E5").Value = "Ready for handoff"
$wsOv.Range("F5").Value = "Ready for handoff"
$wsOv.Range("G5").Value = "2016-11-14 08:01:41"
$wsOv.Hyperlinks.Add($wsOv.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/197298b3fc8f98bf805a29d8aad58b75b0627904/e2e/d1e0eeed-5c24-49dd-bf01-14fb964d9451.md", "", "", "e2e\d1e0eeed-5c24-49dd-bf01-14fb964d9451.md") | Out-Null

$wsOv.Range("A6").Value = "133fdbd1-487f-424b-9b23-57cb4fc069ab.png"
$wsOv.Range("C6").Value = ".png"
$wsOv.Range("E6").Value = "Ready for handoff"
$wsOv.Range("F6").Value = "Ready for handoff"
$wsOv.Range("G6").Value = "2016-11-14 08:01:41"
$wsOv.Hyperlinks.Add($wsOv.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/197298b3fc8f98bf805a29d8aad58b75b0627904/e2e/133fdbd1-487f-424b-9b23-57cb4fc069ab.png", "", "", "e2e\133fdbd1-487f-424b-9b23-57cb4fc069ab.png") | Out-Null

$wsOv.Range("A7").Value = "bb210250-c8aa-4414-af3b-78dac03cd21b.png"
$wsOv.Range("C7").Value = ".png"
$wsOv.Range("E7").Value = "Ready for handoff"
$wsOv.Range("F7").Value = "Ready for handoff"
$wsOv.Range("G7").Value = "2016-11-14 08:01:41"
$wsOv.Hyperlinks.Add($wsOv.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/197298b3fc8f98bf805a29d8aad58b75b0627904/e2e/bb210250-c8aa-4414-af3b-78dac03cd21b.png", "", "", "e2e\bb210250-c8aa-4414-af3b-78dac03cd21b.png") | Out-Null

$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G7"))

# ================= zh-cn sheet =================
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "False"
$wsZh.Range("G5").Value = "d1e0eeed-5c24-49dd-bf01-14fb964d9451.5041514574143389600522da2a6743dec7b832d5.zh-cn.xlf"
$wsZh.Range("H5").Value = "2016-11-14 08:01:26"
$wsZh.Range("K5").Value = "0001-01-01 00:00:00"
$wsZh.Range("M5").Value = "True"
$wsZh.Range("O5").Value = "False"
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/197298b3fc8f98bf805a29d8aad58b75b0627904/e2e/d1e0eeed-5c24-49dd-bf01-14fb964d9451.md", "", "", "d1e0eeed-5c24-49dd-bf01-14fb964d9451.md") | Out-Null

$wsZh.Range("B6").Value = ".png"
$wsZh.Range("C6").Value = "Ready for handoff"
$wsZh.Range("D6").Value = "e2e"
$wsZh.Range("E6").Value = "ht"
$wsZh.Range("F6").Value = "False"
$wsZh.Range("G6").Value = "1f9098e136d87b9ee79d239f202a7f4e9f0c7863.png"
$wsZh.Range("H6").Value = "2016-11-14 08:01:26"
$wsZh.Range("K6").Value = "0001-01-01 00:00:00"
$wsZh.Range("M6").Value = "True(Dependency)"
$wsZh.Range("N6").Value = "e2e\d1e0eeed-5c24-49dd-bf01-14fb964d9451.md"
$wsZh.Range("O6").Value = "False"
$wsZh.Hyperlinks.Add($wsZh.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/197298b3fc8f98bf805a29d8aad58b75b0627904/e2e/133fdbd1-487f-424b-9b23-57cb4fc069ab.png", "", "", "133fdbd1-487f-424b-9b23-57cb4fc069ab.png") | Out-Null

$wsZh.Range("B7").Value = ".png"
$wsZh.Range("C7").Value = "Ready for handoff"
$wsZh.Range("D7").Value = "e2e"
$wsZh.Range("E7").Value = "ht"
$wsZh.Range("F7").Value = "False"
$wsZh.Range("G7").Value = "841abe213388aefc0eaac006053396bd37dac12b.png"
$wsZh.Range("H7").Value = "2016-11-14 08:01:26"
$wsZh.Range("K7").Value = "0001-01-01 00:00:00"
$wsZh.Range("M7").Value = "True(Dependency)"
$wsZh.Range("N7").Value = "e2e\d1e0eeed-5c24-49dd-bf01-14fb964d9451.md"
$wsZh.Range("O7").Value = "False"
$wsZh.Hyperlinks.Add($wsZh.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/197298b3fc8f98bf805a29d8aad58b75b0627904/e2e/bb210250-c8aa-4414-af3b-78dac03cd21b.png", "", "", "bb210250-c8aa-4414-af3b-78dac03cd21b.png") | Out-Null

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P7"))

# ================= de-de sheet =================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "False"
$wsDe.Range("G5").Value = "d1e0eeed-5c24-49dd-bf01-14fb964d9451.5041514574143389600522da2a6743dec7b832d5.de-de.xlf"
$wsDe.Range("H5").Value = "2016-11-14 08:01:41"
$wsDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDe.Range("M5").Value = "True"
$wsDe.Range("O5").Value = "False"
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/197298b3fc8f98bf805a29d8aad58b75b0627904/e2e/d1e0eeed-5c24-49dd-bf01-14fb964d9451.md", "", "", "d1e0eeed-5c24-49dd-bf01-14fb964d9451.md") | Out-Null

$wsDe.Range("B6").Value = ".png"
$wsDe.Range("C6").Value = "Ready for handoff"
$wsDe.Range("D6").Value = "e2e"
$wsDe.Range("E6").Value = "ht"
$wsDe.Range("F6").Value = "False"
$wsDe.Range("G6").Value = "1f9098e136d87b9ee79d239f202a7f4e9f0c7863.png"
$wsDe.Range("H6").Value = "2016-11-14 08:01:41"
$wsDe.Range("K6").Value = "0001-01-01 00:00:00"
$wsDe.Range("M6").Value = "True(Dependency)"
$wsDe.Range("N6").Value = "e2e\d1e0eeed-5c24-49dd-bf01-14fb964d9451.md"
$wsDe.Range("O6").Value = "False"
$wsDe.Hyperlinks.Add($wsDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/197298b3fc8f98bf805a29d8aad58b75b0627904/e2e/133fdbd1-487f-424b-9b23-57cb4fc069ab.png", "", "", "133fdbd1-487f-424b-9b23-57cb4fc069ab.png") | Out-Null

$wsDe.Range("B7").Value = ".png"
$wsDe.Range("C7").Value = "Ready for handoff"
$wsDe.Range("D7").Value = "e2e"
$wsDe.Range("E7").Value = "ht"
$wsDe.Range("F7").Value = "False"
$wsDe.Range("G7").Value = "841abe213388aefc0eaac006053396bd37dac12b.png"
$wsDe.Range("H7").Value = "2016-11-14 08:01:41"
$wsDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDe.Range("M7").Value = "True(Dependency)"
$wsDe.Range("N7").Value = "e2e\d1e0eeed-5c24-49dd-bf01-14fb964d9451.md"
$wsDe.Range("O7").Value = "False"
$wsDe.Hyperlinks.Add($wsDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/197298b3fc8f98bf805a29d8aad58b75b0627904/e2e/bb210250-c8aa-4414-af3b-78dac03cd21b.png", "", "", "bb210250-c8aa-4414-af3b-78dac03cd21b.png") | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P7"))

"done"